$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings as text (some look like plain decimals,
# e.g. "1.00", "130.17"). Force NumberFormat to Text ("@") on each cell
# before writing so Excel does not coerce them into numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.842.32'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.477.46'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '415.87'
$ws.Range("E5").Value = '  +1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.17'
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  -1.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.727'
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("E10").Value = '  +8.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.50'
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("E12").Value = '  +5.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.032.41'
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.57'
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.496.66'
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.67'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.809.45'
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '464.57'
$ws.Range("E21").Value = '  +3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.29'
$ws.Range("E23").Value = '  +2.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.29'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.83'
$ws.Range("E25").Value = '  +15.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.33'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.36'
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.80'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.57'
$ws.Range("E29").Value = '  -2.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.17'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.170'
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.98'
$ws.Range("E34").Value = '  -4.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.15'
$ws.Range("E36").Value = '  +8.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0491'
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.05'
$ws.Range("E39").Value = '  +2.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.74'
$ws.Range("E40").Value = '  +7.38%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.135'
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.75'
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.322'
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.33'
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.41'
$ws.Range("E45").Value = '  +1.80%  '
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0595'
$ws.Range("E47").Value = '  +40.53%  '
$ws.Range("E48").Value = '  +10.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.39'
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.36'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.143'
$ws.Range("E51").Value = '  -2.30%  '
